$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values were regenerated (Strike# -> K), updating rows 2-6.
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 3
$ws.Range("G6").Value = 2
